# Applies the weekly fruit/vegetable price update to the 'Ciruela' sheet.
# Updates date, variety, quality, volume, price and unit/origin columns
# for rows 2-16 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 44243
$ws.Range("K2").Value = 'Black Amber'
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 14500
$ws.Range("S2").Value = 806
$ws.Range("D3").Value = 44245
$ws.Range("M3").Value = 250
$ws.Range("Q3").Value = '$/bandeja 18 kilos granel'
$ws.Range("D4").Value = 44229
$ws.Range("K4").Value = 'Fortuna'
$ws.Range("M4").Value = 300
$ws.Range("D5").Value = 44239
$ws.Range("K5").Value = 'Fortuna'
$ws.Range("L5").Value = 'Primera'
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 15500
$ws.Range("S5").Value = 861
$ws.Range("D6").Value = 44169
$ws.Range("K6").Value = 'Angeleno'
$ws.Range("L6").Value = 'Tercera'
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 24500
$ws.Range("S6").Value = 1361
$ws.Range("D7").Value = 44574
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 19000
$ws.Range("P7").Value = 18500
$ws.Range("S7").Value = 1028
$ws.Range("D8").Value = 44174
$ws.Range("K8").Value = 'Angeleno'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 270
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 21000
$ws.Range("P8").Value = 20500
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 1139
$ws.Range("D9").Value = 44238
$ws.Range("K9").Value = 'Black Amber'
$ws.Range("L9").Value = 'Segunda'
$ws.Range("D10").Value = 44238
$ws.Range("K10").Value = 'Fortuna'
$ws.Range("M10").Value = 300
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 14500
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 806
$ws.Range("D11").Value = 44580
$ws.Range("K11").Value = 'Black Amber'
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 270
$ws.Range("N11").Value = 19000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 19500
$ws.Range("R11").Value = 'Región Metropolitana'
$ws.Range("S11").Value = 1083
$ws.Range("D12").Value = 44217
$ws.Range("L12").Value = 'Segunda'
$ws.Range("N12").Value = 16000
$ws.Range("O12").Value = 17000
$ws.Range("P12").Value = 16500
$ws.Range("R12").Value = 'Región Metropolitana'
$ws.Range("S12").Value = 917
$ws.Range("D13").Value = 44175
$ws.Range("K13").Value = 'Angeleno'
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 21000
$ws.Range("O13").Value = 22000
$ws.Range("P13").Value = 21500
$ws.Range("S13").Value = 1194
$ws.Range("D14").Value = 44285
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 14000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 14500
$ws.Range("S14").Value = 806
$ws.Range("D15").Value = 44278
$ws.Range("K15").Value = 'Angeleno'
$ws.Range("L15").Value = 'Primera'
$ws.Range("N15").Value = 15000
$ws.Range("O15").Value = 16000
$ws.Range("P15").Value = 15500
$ws.Range("Q15").Value = '$/caja 18 kilos granel'
$ws.Range("S15").Value = 861
$ws.Range("D16").Value = 44314
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 250
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 14500
$ws.Range("Q16").Value = '$/bandeja 18 kilos granel'
$ws.Range("S16").Value = 806
